# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stocks) sheet gets three new trailing columns appended to its
# existing table (name, owner, quantity, face_value, currency, total):
#   H: date             -> "2012-05-01"   (filing date for this report)
#   I: legislator_name  -> "翁重鈞"       (the legislator this report belongs to)
#   J: legislator_id    -> 551            (legislator's numeric id)
#
# Every data row gets the same date / legislator_name / legislator_id, since
# they describe the report as a whole rather than the individual stock entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorName = "翁重鈞"
$legislatorId = 551
$reportDate = "2012-05-01"

$lastRow = 4

# Match the look of the existing table: copy the header cell's formatting
# (bold, centered, bordered) onto the new header cells, and the existing
# body-row formatting onto the new body cells, before filling in values.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

$ws.Range("G2:G" + $lastRow).Copy()
$ws.Range("H2:J" + $lastRow).PasteSpecial(-4122)

# Header row
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# Force the date column to be stored as plain text (not auto-converted to a
# date serial number) so it round-trips as the literal string "2012-05-01".
$ws.Range("H2:H" + $lastRow).NumberFormat = "@"

# Data rows -> same date / legislator_name / legislator_id for every row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $reportDate
    $ws.Cells.Item($r, 9).Value = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}
